# Commit: "added classification by random forest"
#
# The evaluation was re-run after adding a Random Forest (RF) classifier to the
# pipeline. This changed the metric values for every existing model (kNN, SVM,
# LR, NB) because of the new run, inserted a brand-new "RF" results row, and
# pushed the "Ensemble" row (which now also folds in RF) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row labels (A6 = "RF", A7 = "Ensemble") -------------------
# Copy the existing header-style formatting (bold, bordered, centered) from A5
# onto the two new label cells before setting their text, so the new cells pick
# up the same look as the rest of column A rather than the default style.
$ws.Range("A5").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)

$ws.Range("A6").Value = "RF"
$ws.Range("A7").Value = "Ensemble"

# --- Write the (re-computed) metric values for every row --------------------

# Row 2: kNN
$row2 = New-Object 'object[,]' 1,24
$row2[0,0] = 0.6138183482040723
$row2[0,1] = 0.694119958807532
$row2[0,2] = 0.6138183482040723
$row2[0,3] = 0.6109102501821633
$row2[0,4] = 0.6760009151223976
$row2[0,5] = 0.7340059640652733
$row2[0,6] = 0.6760009151223976
$row2[0,7] = 0.677086542942525
$row2[0,8] = 0.6759322809425761
$row2[0,9] = 0.7109154473037099
$row2[0,10] = 0.6759322809425761
$row2[0,11] = 0.6754991818369545
$row2[0,12] = 0.7917867764813545
$row2[0,13] = 0.8077005426910869
$row2[0,14] = 0.7917867764813545
$row2[0,15] = 0.791567739767727
$row2[0,16] = 0.7061541981239993
$row2[0,17] = 0.7767767321584371
$row2[0,18] = 0.7061541981239993
$row2[0,19] = 0.7108778336318742
$row2[0,20] = 0.6889727751086708
$row2[0,21] = 0.7929746933734337
$row2[0,22] = 0.6889727751086708
$row2[0,23] = 0.6977995982724715
$ws.Range("B2:Y2").Value = $row2

# Row 3: SVM
$row3 = New-Object 'object[,]' 1,24
$row3[0,0] = 0.8004346831388698
$row3[0,1] = 0.8054812875657584
$row3[0,2] = 0.8004346831388698
$row3[0,3] = 0.8008134154112974
$row3[0,4] = 0.8070006863417982
$row3[0,5] = 0.81867355006829
$row3[0,6] = 0.8070006863417982
$row3[0,7] = 0.8089305373102448
$row3[0,8] = 0.7983070235644018
$row3[0,9] = 0.8074781483761285
$row3[0,10] = 0.7983070235644018
$row3[0,11] = 0.7993243619866995
$row3[0,12] = 0.8304049416609471
$row3[0,13] = 0.8353814862452783
$row3[0,14] = 0.8304049416609471
$row3[0,15] = 0.8307171212580075
$row3[0,16] = 0.8476321207961565
$row3[0,17] = 0.853558973059082
$row3[0,18] = 0.8476321207961565
$row3[0,19] = 0.8475142028388465
$row3[0,20] = 0.8476549988560971
$row3[0,21] = 0.8534647084786318
$row3[0,22] = 0.8476549988560971
$row3[0,23] = 0.8477463511861221
$ws.Range("B3:Y3").Value = $row3

# Row 4: LR
$row4 = New-Object 'object[,]' 1,24
$row4[0,0] = 0.8412033859528713
$row4[0,1] = 0.8446366039755485
$row4[0,2] = 0.8412033859528713
$row4[0,3] = 0.8407444217568611
$row4[0,4] = 0.8454358270418668
$row4[0,5] = 0.8505367336269611
$row4[0,6] = 0.8454358270418668
$row4[0,7] = 0.8452030091656443
$row4[0,8] = 0.8411805078929306
$row4[0,9] = 0.844549349730735
$row4[0,10] = 0.8411805078929306
$row4[0,11] = 0.8413400401085502
$row4[0,12] = 0.8305193319606496
$row4[0,13] = 0.8368399184355193
$row4[0,14] = 0.8305193319606496
$row4[0,15] = 0.8296821578910329
$row4[0,16] = 0.8369709448638755
$row4[0,17] = 0.8423092278583123
$row4[0,18] = 0.8369709448638755
$row4[0,19] = 0.8366381363445029
$row4[0,20] = 0.8391443605582246
$row4[0,21] = 0.8451562018348253
$row4[0,22] = 0.8391443605582246
$row4[0,23] = 0.8389726415160503
$ws.Range("B4:Y4").Value = $row4

# Row 5: NB
$row5 = New-Object 'object[,]' 1,24
$row5[0,0] = 0.8132235186456189
$row5[0,1] = 0.8215943879052988
$row5[0,2] = 0.8132235186456189
$row5[0,3] = 0.8130234233276126
$row5[0,4] = 0.8454815831617479
$row5[0,5] = 0.8517901636782993
$row5[0,6] = 0.8454815831617479
$row5[0,7] = 0.8459302328161922
$row5[0,8] = 0
$row5[0,9] = 0
$row5[0,10] = 0
$row5[0,11] = 0
$row5[0,12] = 0.8346831388698238
$row5[0,13] = 0.8474988865737526
$row5[0,14] = 0.8346831388698238
$row5[0,15] = 0.8325940159128592
$row5[0,16] = 0.847609242736216
$row5[0,17] = 0.8565718259019214
$row5[0,18] = 0.847609242736216
$row5[0,19] = 0.8450343347222546
$row5[0,20] = 0
$row5[0,21] = 0
$row5[0,22] = 0
$row5[0,23] = 0
$ws.Range("B5:Y5").Value = $row5

# Row 6: RF
$row6 = New-Object 'object[,]' 1,24
$row6[0,0] = 0.8155113246396706
$row6[0,1] = 0.8241364953466619
$row6[0,2] = 0.8155113246396706
$row6[0,3] = 0.8153985638556321
$row6[0,4] = 0.8198581560283689
$row6[0,5] = 0.8290745711167341
$row6[0,6] = 0.8198581560283689
$row6[0,7] = 0.8207637069532379
$row6[0,8] = 0.7447494852436515
$row6[0,9] = 0.7512047604522312
$row6[0,10] = 0.7447494852436515
$row6[0,11] = 0.743942986068282
$row6[0,12] = 0.7961793639899337
$row6[0,13] = 0.806091831194253
$row6[0,14] = 0.7961793639899337
$row6[0,15] = 0.7965873795539935
$row6[0,16] = 0.7983527796842829
$row6[0,17] = 0.8069839088129935
$row6[0,18] = 0.7983527796842829
$row6[0,19] = 0.7989084405626932
$row6[0,20] = 0.8091054678563259
$row6[0,21] = 0.8167073583053162
$row6[0,22] = 0.8091054678563259
$row6[0,23] = 0.8084391210551228
$ws.Range("B6:Y6").Value = $row6

# Row 7: Ensemble
$row7 = New-Object 'object[,]' 1,24
$row7[0,0] = 0.8390986044383437
$row7[0,1] = 0.8438015173036506
$row7[0,2] = 0.8390986044383437
$row7[0,3] = 0.839706003705081
$row7[0,4] = 0.8670327156257149
$row7[0,5] = 0.8689161636703311
$row7[0,6] = 0.8670327156257149
$row7[0,7] = 0.8661392521523243
$row7[0,8] = 0.8262411347517731
$row7[0,9] = 0.8397427235955638
$row7[0,10] = 0.8262411347517731
$row7[0,11] = 0.8273741760395279
$row7[0,12] = 0.8498055364905056
$row7[0,13] = 0.8547738962830355
$row7[0,14] = 0.8498055364905056
$row7[0,15] = 0.8501608646320928
$row7[0,16] = 0.87124227865477
$row7[0,17] = 0.8768968005363064
$row7[0,18] = 0.87124227865477
$row7[0,19] = 0.8709454940034813
$row7[0,20] = 0.8411805078929306
$row7[0,21] = 0.8466553609676971
$row7[0,22] = 0.8411805078929306
$row7[0,23] = 0.8416830831706434
$ws.Range("B7:Y7").Value = $row7

Write-Host "Applied RF row insertion + metric refresh."
